# Update the "Förändrad" (column C) date for every existing data row
# (rows 2-156) from 2023-09-23 (45192) to 2023-10-03 (45202).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C156").Value = 45202

# Row 156 becomes an explicit 15pt row (ht="15" customHeight="1"), matching
# the new rows being appended below it.
$ws.Rows.Item(156).RowHeight = 15

# Append the 7 new cleaning notifications (rows 157-163).
$newRows = @(
    @{ Row = 157; A = "A 47013-2023"; B = 45196; C = 45202; D = "STOCKHOLMS LÄN"; E = "NYNÄSHAMN"; G = 2.3 },
    @{ Row = 158; A = "A 47001-2023"; B = 45196; C = 45202; D = "STOCKHOLMS LÄN"; E = "NYNÄSHAMN"; G = 1.1 },
    @{ Row = 159; A = "A 47098-2023"; B = 45196; C = 45202; D = "STOCKHOLMS LÄN"; E = "NYNÄSHAMN"; G = 8 },
    @{ Row = 160; A = "A 47115-2023"; B = 45196; C = 45202; D = "STOCKHOLMS LÄN"; E = "NYNÄSHAMN"; G = 1.7 },
    @{ Row = 161; A = "A 47093-2023"; B = 45196; C = 45202; D = "STOCKHOLMS LÄN"; E = "NYNÄSHAMN"; G = 1 },
    @{ Row = 162; A = "A 47110-2023"; B = 45196; C = 45202; D = "STOCKHOLMS LÄN"; E = "NYNÄSHAMN"; G = 0.9 },
    @{ Row = 163; A = "A 47089-2023"; B = 45196; C = 45202; D = "STOCKHOLMS LÄN"; E = "NYNÄSHAMN"; G = 0.8 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.A

    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 2).NumberFormat = "YYYY-MM-DD"

    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 3).NumberFormat = "YYYY-MM-DD"

    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E

    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = 0
    $ws.Cells.Item($row, 9).Value = 0
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = 0
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0
    $ws.Cells.Item($row, 15).Value = 0
    $ws.Cells.Item($row, 16).Value = 0
    $ws.Cells.Item($row, 17).Value = 0

    $ws.Cells.Item($row, 18).Value = ""
    $ws.Cells.Item($row, 18).WrapText = $true

    # Rows 157-162 get an explicit 15pt row height; the final row (163) does
    # not, matching the source workbook.
    if ($row -ne 163) {
        $ws.Rows.Item($row).RowHeight = 15
    }
}
